$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows("2:2").Delete() | Out-Null
$ws.Rows("2:2").Select() | Out-Null
